# Update "想去人数" (want-to-go count) values in the F column across sheets,
# reflecting a refreshed data scrape (per commit "Update gh-pages to output
# generated at 456a3b4").

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F6").Value = 2500
$wsExhibit.Range("F12").Value = 33
$wsExhibit.Range("F13").Value = 7237
$wsExhibit.Range("F14").Value = 319
$wsExhibit.Range("F20").Value = 8630
$wsExhibit.Range("F38").Value = 1176
$wsExhibit.Range("F41").Value = 3727
$wsExhibit.Range("F46").Value = 201

# Sheet "演出" (Performance)
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F16").Value = 22

# Sheet "全部类型" (All Types, aggregate of all the above)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F7").Value = 2500
$wsAll.Range("F16").Value = 33
$wsAll.Range("F17").Value = 7237
$wsAll.Range("F18").Value = 319
$wsAll.Range("F23").Value = 8630
$wsAll.Range("F39").Value = 1176
$wsAll.Range("F42").Value = 3727
$wsAll.Range("F45").Value = 22
$wsAll.Range("F47").Value = 201
